# Auto-generated edit script applying the cryptos.xlsx diff
# Forces each updated cell to remain plain TEXT (matching the original
# inlineStr cell type) by prefixing the assigned value with a leading
# apostrophe (Excel's literal-text marker) and then resetting the cell
# style back to "Normal" so no stray number format / quote-prefix style
# sticks to the cell (keeps output identical in shape to the source).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.431.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +6.82%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.112.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +4.33%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'586.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.59%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'145.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.14%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.107.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +4.36%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +2.40%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +14.48%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +7.35%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +3.95%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +8.48%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'35.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +5.40%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.74%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.625.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.29%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'7.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.14%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'63.284.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +6.52%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.107.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.15%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'467.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.93%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.727"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.86%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +7.48%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'13.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.43%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'81.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.46%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +11.59%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.06%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +4.81%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.11%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +9.49%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'27.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.76%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +5.31%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0₃0878"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +13.49%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +16.62%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +6.76%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +21.47%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.39%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'50.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.07%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'440.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +9.90%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'8.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'2.914.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.58%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0369"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.96%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +11.48%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +5.08%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +8.44%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'35.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.59%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D49").Value = "'123.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.11%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +1.47%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'24.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +6.09%  "
$ws.Range("E51").Style = "Normal"
